# Updated hyperlink, shave 50% off the per-board price by changing 2-pin terminal block
# suppliers (TE Connectivity -> DBParts/Amazon) in the LED Matrix PCB BOM.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: swap the 2-pin terminal block part for the cheaper Amazon/DBParts one ---
$newName = "DBParts 20 pc 2-pin 0.1” pitch Terminal Blocks"

$ws.Range("A5").Value = $newName
$ws.Range("B5").Value = $newName
$ws.Range("C5").Value = 6.99
$ws.Range("D5").Formula = "=73/20"
$ws.Range("F5").Value = "Amazon"
$ws.Range("G5").Value = "https://www.amazon.com/DBParts-20pcs-Terminal-Connector-2-54mm/dp/B07NSJV6NW/ref=sxbs_sxwds-stvp?cv_ct_cx=terminal+block+assortment&keywords=terminal+block+assortment&pd_rd_i=B07NSJV6NW&pd_rd_r=f9187fb7-4ab9-45b4-a445-3bf2b68a1d13&pd_rd_w=WAjjn&pd_rd_wg=olYB9&pf_rd_p=a6d018ad-f20b-46c9-8920-433972c7d9b7&pf_rd_r=PMA44C3EHR468DBRJFP7&psc=1&qid=1581446551&sr=1-3-dd5817a1-1ba7-46c2-8996-f96e7b0f409c"

# Give the link cell the same "link-styled" blue text used elsewhere in the sheet
$ws.Range("G5").Font.Color = 16711680

# --- Column widths (A and B got wider) ---
$ws.Columns.Item(1).ColumnWidth = 38.833333333
$ws.Columns.Item(2).ColumnWidth = 40.666666666

# --- View: zoom in to 120% and move the selection to B15 ---
$excel.ActiveWindow.Zoom = 120
$excel.ActiveWindow.DisplayGridlines = $true
$ws.Range("B15").Select() | Out-Null
